# Update column G ("K") values on the active worksheet.
# These are the Strike# -> K replacement values computed by the
# save_data regeneration described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 1
    4  = 6
    5  = 10
    6  = 3
    7  = 8
    8  = 6
    9  = 7
    10 = 4
    11 = 6
    12 = 4
    13 = 7
    14 = 5
    15 = 9
    16 = 8
    17 = 5
    18 = 7
    19 = 5
    20 = 8
    21 = 5
    22 = 7
    23 = 3
    24 = 2
    25 = 6
    26 = 0
    27 = 9
    28 = 2
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
